# Ad hoc fix to pull_data
# pull_data is unable to retrieve annual data at the moment, so the date
# column is rewritten from Excel date-serial numbers (formatted yyyy-mm-dd)
# to plain "YYYY-12-31" text, and the forecast-only usphpi/casusxam values
# for 2020-2030 (rows 52-62) are removed since they can no longer be pulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 62

# Column A currently holds date serials with a custom "yyyy-mm-dd" number
# format (style index 2 in the original file). Writing a date-look-alike
# string to .Value would just get reinterpreted as that same date serial,
# so force literal text with a leading apostrophe, then drop the range
# back to the default "Normal" style so no per-cell date-style reference
# remains on these cells (matches the original, unstyled A column).
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $year = $row + 1968
    $ws.Cells.Item($row, 1).Value = "'$year-12-31"
}

$ws.Range("A$firstRow`:A$lastRow").Style = "Normal"

# usphpi (B) and casusxam (C) for 2020-2030 (rows 52-62) can no longer be
# retrieved, so clear those cells entirely rather than leaving stale data.
$ws.Range("B52:C62").ClearContents()
